# Fruta / hortaliza, semanal
# Insert the week's new Chirimoya price records (Femacal de La Calera,
# Provincia de Limarí) at the top of the data and push the existing
# rows down by two positions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at row 7 (pushes old rows 7..30 down to 9..32,
# inheriting the date number-format from the row being displaced).
$ws.Rows.Item(7).Insert()
$ws.Rows.Item(7).Insert()

# New row 7: Primera
$ws.Range("A7").Value = 3
$ws.Range("B7").Value = "Femacal de La Calera"
$ws.Range("C7").Value = "Coquimbo"
$ws.Range("D7").Value = 44459
$ws.Range("E7").Value = 5
$ws.Range("F7").Value = "Fruta"
$ws.Range("G7").Value = 100107
$ws.Range("H7").Value = "Otros"
$ws.Range("I7").Value = 100107002
$ws.Range("J7").Value = "Chirimoya"
$ws.Range("K7").Value = "Cultivar IV Región"
$ws.Range("L7").Value = "Primera"
$ws.Range("M7").Value = 75
$ws.Range("N7").Value = 26500
$ws.Range("O7").Value = 27000
$ws.Range("P7").Value = 26767
$ws.Range("Q7").Value = "$/bandeja 10 kilos"
$ws.Range("R7").Value = "Provincia de Limarí"
$ws.Range("S7").Value = 2677
$ws.Range("T7").Value = 10

# New row 8: Segunda
$ws.Range("A8").Value = 3
$ws.Range("B8").Value = "Femacal de La Calera"
$ws.Range("C8").Value = "Coquimbo"
$ws.Range("D8").Value = 44459
$ws.Range("E8").Value = 5
$ws.Range("F8").Value = "Fruta"
$ws.Range("G8").Value = 100107
$ws.Range("H8").Value = "Otros"
$ws.Range("I8").Value = 100107002
$ws.Range("J8").Value = "Chirimoya"
$ws.Range("K8").Value = "Cultivar IV Región"
$ws.Range("L8").Value = "Segunda"
$ws.Range("M8").Value = 40
$ws.Range("N8").Value = 24000
$ws.Range("O8").Value = 24000
$ws.Range("P8").Value = 24000
$ws.Range("Q8").Value = "$/bandeja 10 kilos"
$ws.Range("R8").Value = "Provincia de Limarí"
$ws.Range("S8").Value = 2400
$ws.Range("T8").Value = 10
